# Updated symbol list on Mon Dec 12 09:21:30 UTC 2022 with GitHub Actions
#
# The "Price" column (D) holds text-looking numeric strings (e.g. "281.55").
# Writing a numeric-looking string straight into Range.Value makes Excel
# coerce it to a real number, which would change the cell's stored type
# (string -> number) and not match the source data (which keeps these as
# text). To force Excel to keep the literal text, we briefly mark the cell
# as Text (NumberFormat "@") before assigning, then clear the format again
# so the cell's style index is left exactly as it was (General / default).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.ClearFormats()
}

Set-TextValue "D2"  "281.72"
Set-TextValue "D3"  "20.87"
Set-TextValue "D4"  "6.236"
Set-TextValue "D5"  "0.06147"
Set-TextValue "D6"  "3.573"
Set-TextValue "D7"  "6.558"
Set-TextValue "D8"  "1.483"
Set-TextValue "D9"  "0.8171"
Set-TextValue "D10" "0.01383"
Set-TextValue "D11" "0.1628"
Set-TextValue "D12" "0.08284"
Set-TextValue "D13" "0.03551"
Set-TextValue "D14" "0.03189"
Set-TextValue "D15" "0.09138"
Set-TextValue "D16" "3.721"
Set-TextValue "D18" "0.04641"
Set-TextValue "D19" "0.006420"
Set-TextValue "D20" "0.006178"
Set-TextValue "D23" "3.808"
Set-TextValue "D24" "2.337"
Set-TextValue "D25" "0.3373"
Set-TextValue "D40" "0.04669"
Set-TextValue "D41" "0.007101"
Set-TextValue "D42" "0.1102"
Set-TextValue "D43" "0.003509"
Set-TextValue "D44" "0.01133"
Set-TextValue "D45" "0.00006314"
Set-TextValue "D47" "0.9997"
Set-TextValue "D48" "0.002941"
Set-TextValue "D49" "0.00001899"
